$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 70: add EARNED value ---
$ws.Range("C70").Value = 1.25

# --- Row 71: add PERIOD date ---
$ws.Range("A71").Value = 44774

# --- Row 72: add PERIOD date + EARNED value ---
$ws.Range("A72").Value = 44805
$ws.Range("C72").Value = 1.25

# --- Row 73: add PERIOD date + EARNED value ---
$ws.Range("A73").Value = 44835
$ws.Range("C73").Value = 1.25

# --- Row 74: add PERIOD date + EARNED value ---
$ws.Range("A74").Value = 44866
$ws.Range("C74").Value = 1.25

# --- Row 76: becomes the "2023" year-header row (mirrors rows 10/24/37/50/63) ---
# (set before row 75's text cells so the new shared-string order matches: 2023, SL(2-0-0), 12/9,15/2022)
$cellA76 = $ws.Range("A76")
$cellA76.NumberFormat = "@"
$cellA76.Value = "2023"
$ws.Range("A63").Copy() | Out-Null
$cellA76.PasteSpecial(-4122) | Out-Null

# --- Row 75: add PERIOD date, PARTICULARS, EARNED, Absence Undertime W/ Pay (2nd set), REMARKS ---
$ws.Range("A75").Value = 44896
$ws.Range("B75").Value = "SL(2-0-0)"
$ws.Range("C75").Value = 1.25
$ws.Range("H75").Value = 2
$ws.Range("K75").Value = "12/9,15/2022"

# --- Row 77: add PERIOD date + EARNED value ---
$ws.Range("A77").Value = 44927
$ws.Range("C77").Value = 1.25

# --- Row 78: add PERIOD date + EARNED value ---
$ws.Range("A78").Value = 44958
$ws.Range("C78").Value = 1.25

# --- Remove the now-unused blank row 129 (old row 130 shifts up to become the new last row 129) ---
$ws.Rows.Item(129).Delete()

$excel.CutCopyMode = 0
